# Migration of the "simple" metadata namespace to "datamodel" on the
# DataModels sheet, plus refresh of the saved selection/cursor position,
# matching the "Complete migration with tests working" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataModels")

# Rename the custom-metadata key headers from the old "excel.simple"
# namespace to the new "excel.datamodel" namespace.
$ws.Range("G1").Value = "uk.ac.ox.softeng.maurodatamapper.plugins.excel.datamodel:reviewed"
$ws.Range("H1").Value = "uk.ac.ox.softeng.maurodatamapper.plugins.excel.datamodel:approved"
$ws.Range("I1").Value = "uk.ac.ox.softeng.maurodatamapper.plugins.excel.datamodel:distributed"

# Move the cursor/selection on the DataModels sheet to reflect the
# latest author session (was H3, now J3).
$ws.Activate()
$ws.Range("J3").Select()
